# Remove some old notes.
#
# Slide 5 speaker notes: drop the "I'll discuss two pieces of work..." line.
# Slide 6 speaker notes: drop the "Skip if required." line.
# In both cases the notes placeholder ends up with no text (an empty
# paragraph), matching the diff's endParaRPr-only paragraph.

$p = $ppt.ActivePresentation

$slide5 = $p.Slides.Item(5)
$notes5 = $slide5.NotesPage
$notesBody5 = $notes5.Shapes.Item(2)
$notesBody5.TextFrame.TextRange.Text = ""

$slide6 = $p.Slides.Item(6)
$notes6 = $slide6.NotesPage
$notesBody6 = $notes6.Shapes.Item(2)
$notesBody6.TextFrame.TextRange.Text = ""
